$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '27.390.14'
$c.Style = 'Normal'
$ws.Range('E2').Value = '  -1.75%  '

$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '1.656.49'
$c.Style = 'Normal'
$ws.Range('E3').Value = '  -0.48%  '

$ws.Range('E4').Value = '  -0.07%  '

$ws.Range('E5').Value = '  -0.65%  '

$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '0.514'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -0.43%  '

$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  -0.10%  '

$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '23.62'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  +0.49%  '

$ws.Range('E9').Value = '  -0.44%  '

$ws.Range('E10').Value = '  -1.15%  '

$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.0876'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -0.41%  '

$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '1.890.40'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -0.43%  '

$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '1.651.96'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  -0.70%  '

$ws.Range('E14').Value = '  -1.61%  '

$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '0.568'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  +3.34%  '

$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '65.66'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -0.57%  '

$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '27.379.54'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  -1.55%  '

$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '231.59'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -6.96%  '

$ws.Range('E19').Value = '  -0.75%  '

$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '7.48'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -0.14%  '

$ws.Range('E21').Value = '  -0.06%  '

$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '4.37'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -2.25%  '

$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '9.35'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +2.68%  '

$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '2.02'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -1.21%  '

$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '147.71'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +0.63%  '

$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '7.11'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -1.18%  '

$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '15.90'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -2.30%  '

$ws.Range('E28').Value = '  +0.02%  '

$ws.Range('E29').Value = '  -0.42%  '

$ws.Range('E30').Value = '  -0.83%  '

$ws.Range('E31').Value = '  -4.44%  '

$ws.Range('E32').Value = '  -1.49%  '

$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '1.435.13'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +0.04%  '

$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '3.14'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +0.18%  '

$ws.Range('E35').Value = '  +0.38%  '

$ws.Range('E36').Value = '  -0.74%  '

$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.908'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -2.38%  '

$ws.Range('E38').Value = '  -1.85%  '

$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.0169'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +0.18%  '

$ws.Range('E40').Value = '  -0.10%  '

$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -0.12%  '

$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '5.51'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +1.91%  '

$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '64.97'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -6.80%  '

$ws.Range('E45').Value = '  +0.17%  '

$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '1.798.47'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -0.39%  '

$ws.Range('E47').Value = '  -1.02%  '

$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '87.99'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -1.42%  '

$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.0₆0106'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -1.37%  '

$ws.Range('E50').Value = '  -0.35%  '

$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '7.74'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -0.64%  '
